$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 766
$ws1.Range("F3").Value = 606
$ws1.Range("F5").Value = 4038
$ws1.Range("F7").Value = 8585
$ws1.Range("F8").Value = 221
$ws1.Range("F9").Value = 495
$ws1.Range("F10").Value = 67
$ws1.Range("F11").Value = 524

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 18
$ws2.Range("F4").Value = 4
$ws2.Range("F5").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 766
$ws4.Range("F3").Value = 606
$ws4.Range("F7").Value = 4038
$ws4.Range("F9").Value = 4
$ws4.Range("F10").Value = 8585
$ws4.Range("F11").Value = 221
$ws4.Range("F12").Value = 495
$ws4.Range("F13").Value = 67
$ws4.Range("F15").Value = 1
$ws4.Range("F16").Value = 524
